# Actualización automática hashcode vie jul 26 01:53:18 CEST 2019
# Updates the "hashcode" value (column B) for the rows identified by their
# code (column A) in the active worksheet, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Code = "05-050301A";   New = "1f682c4baf00039722b9d3b2a8f6431f" },
    @{ Code = "05-050105A";   New = "775da89266fde57dfe7ca7c89abf5d91" },
    @{ Code = "05-050103A";   New = "8a74666dc4ebb183229cedc771aa374f" },
    @{ Code = "05-050104A";   New = "e5a9c26e094a5557ae9c4aa83e416d55" },
    @{ Code = "05-050101A";   New = "0c473cacc596f7b80f753639d0d0ca9c" },
    @{ Code = "05-050102A";   New = "8c9098805d070995ea6995c660cc73a1" },
    @{ Code = "05-050301TP";  New = "81667d4f5140992663fc6287a415e11f" },
    @{ Code = "05-050007TC";  New = "0164192226833e8b2508d9634b0ba903" },
    @{ Code = "05-050007TP";  New = "adf3c1215f1ec05392a34e4fcab6d818" },
    @{ Code = "05-050105TC";  New = "1e5c3f3bf56fea72588394470e1cc359" },
    @{ Code = "05-050105TP";  New = "bc95cae257a5ff8399d8aa38ac0096e0" },
    @{ Code = "05-050101TP";  New = "9283cf6e227051ed64790cd8214746ac" },
    @{ Code = "05-050103TP";  New = "3d3502f758d76be92c0f4e2ea3201dd1" },
    @{ Code = "05-050006A";   New = "b4d216af1c0225064ccc574065e16246" },
    @{ Code = "05-050201A";   New = "61c4f18193adac7d146bc75c0f680430" },
    @{ Code = "05-050007A";   New = "8317bc5e1079993b6d686cc7d773b4ef" },
    @{ Code = "05-050102TP";  New = "856d009b685edcaa25e7aebd1e4cb92c" },
    @{ Code = "05-050006TC";  New = "5f1e48ea2ee37ac4a0cd6534daf28e1d" },
    @{ Code = "05-050006TP";  New = "deeeabb02d47e448e34e5d3bbaeb8dad" },
    @{ Code = "05-050104TC";  New = "831b12f239db1883cfb6a62cd480eabe" },
    @{ Code = "05-050104TM";  New = "e0b748b7abab51601ff88878e1646e1d" },
    @{ Code = "05-050104TP";  New = "e72e4ad52475855fd285dd2b5bbecbd4" }
)

foreach ($u in $updates) {
    $found = $ws.Columns.Item(1).Find($u.Code)
    if ($found -eq $null) {
        Write-Output ("NOT FOUND: " + $u.Code)
    } else {
        $row = $found.Row
        $ws.Cells.Item($row, 2).Value = $u.New
    }
}
